$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.308.71"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "2.184.66"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'255.31"
$ws.Range("E5").Value = "  +5.53%  "

$ws.Range("E6").Value = "  +0.78%  "

$ws.Range("D7").Value = "'68.04"
$ws.Range("E7").Value = "  -1.69%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.575"
$ws.Range("E9").Value = "  +6.00%  "

$ws.Range("D10").Value = "'37.59"
$ws.Range("E10").Value = "  +3.33%  "

$ws.Range("D11").Value = "'59.07"
$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("D12").Value = "'0.0934"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("D13").Value = "'7.11"
$ws.Range("E13").Value = "  +7.96%  "

$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").Value = "2.502.41"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("D16").Value = "'0.870"
$ws.Range("E16").Value = "  +4.69%  "

$ws.Range("D17").Value = "'14.44"
$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("D18").Value = "2.210.77"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").Value = "41.192.99"
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  +1.50%  "

$ws.Range("E21").Value = "  +2.07%  "

$ws.Range("D22").Value = "'71.78"
$ws.Range("E22").Value = "  -0.96%  "

$ws.Range("D23").Value = "'231.98"
$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("E25").Value = "  +10.02%  "

$ws.Range("D26").Value = "'11.68"
$ws.Range("E26").Value = "  +21.02%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  +5.82%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "'168.80"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").Value = "'20.65"
$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").Value = "'0.0750"
$ws.Range("E33").Value = "  +7.35%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "'5.46"
$ws.Range("E35").Value = "  +6.57%  "

$ws.Range("D36").Value = "'26.45"
$ws.Range("E36").Value = "  +11.05%  "

$ws.Range("D37").Value = "'4.62"
$ws.Range("E37").Value = "  +1.11%  "

$ws.Range("E38").Value = "  +8.07%  "

$ws.Range("D39").Value = "'0.0300"
$ws.Range("E39").Value = "  +11.04%  "

$ws.Range("E40").Value = "  -2.92%  "

$ws.Range("D41").Value = "'12.55"
$ws.Range("E41").Value = "  +20.58%  "

$ws.Range("D42").Value = "'5.66"
$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'5.15"
$ws.Range("E43").Value = "  +5.98%  "

$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "'64.15"
$ws.Range("E44").Value = "  +2.71%  "

$ws.Range("D45").Value = "'0.200"
$ws.Range("E45").Value = "  +4.81%  "

$ws.Range("D46").Value = "'8.63"
$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("E47").Value = "  +2.83%  "

$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("E49").Value = "  +4.70%  "

$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("D51").Value = "'4.25"
$ws.Range("E51").Value = "  -5.27%  "

